# Add a new row (row 24) to the lookup table for the "ConvertorIn" activity,
# following directly after the last existing row (row 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the row above (A:D) so the new row matches the
# table's existing look (fonts, borders, number formats, etc.).
$ws.Range("A23:D23").Copy($ws.Range("A24:D24"))
# Column E in the most recently added rows uses a slightly different style
# than the rest of the table; mirror an earlier row's E-cell formatting
# instead so the new row's IsExist cell matches the intended look.
$ws.Range("E2").Copy($ws.Range("E24"))

# Populate the new row's values.
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "ConvertorIn"
$ws.Range("C24").Value = "/home/pmuser01/converters/convert_biller/in"
$ws.Range("D24").Value = "Yes"
